# Update with latest cht-conf changes again and remove NO_LABEL
#
# The "survey" sheet's row 3 (the "begin_group"/"page" header row) had a
# label of "NO_LABEL" in column C. That value is no longer wanted, so the
# cell is cleared (this also drops "NO_LABEL" from the shared string
# table once the file is re-saved). Column D ("field-list" / appearance)
# is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Range("C3").ClearContents()

# The sheet's frozen bottom-right pane had its active cell parked on A5;
# move it back to A2 (the first data row under the frozen header).
$ws.Activate()
$ws.Range("A2").Select()
